$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(32, 35, 36, 43, 44)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 8).Value = -1
}
